# Auto-generated: refresh market-price snapshot columns (H:N) across all eight
# crafting-class Leve tables, per the scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws.Range("H53").Value = 472.2
$ws.Range("I53").Value = 461
$ws.Range("J53").Value = 475
$ws.Range("K53").Value = 461
$ws.Range("L53").Value = 475
$ws.Range("M53").Value = 176
$ws.Range("N53").Value = -1749

$ws.Range("H69").Value = 56684.7
$ws.Range("I69").Value = 4257.3335
$ws.Range("K69").Value = 12772.0005
$ws.Range("M69").Value = -11898.0005

$ws.Range("H72").Value = 56684.7
$ws.Range("I72").Value = 4257.3335
$ws.Range("K72").Value = 38316.0015
$ws.Range("M72").Value = -33948.0015

$ws.Range("H74").Value = 10900.5
$ws.Range("I74").Value = 10900.5
$ws.Range("K74").Value = 10900.5
$ws.Range("M74").Value = -9964.5

$ws.Range("H77").Value = 10900.5
$ws.Range("I77").Value = 10900.5
$ws.Range("K77").Value = 54502.5
$ws.Range("M77").Value = -49822.5

$ws.Range("H88").Value = 4500
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 8000
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 8000
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -8812

$ws.Range("H91").Value = 4500
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 8000
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 8000
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -10808

$ws.Range("H96").Value = 426.23077
$ws.Range("I96").Value = 299.18182
$ws.Range("K96").Value = 897.54546
$ws.Range("M96").Value = 475.45454

$ws.Range("H99").Value = 823.5
$ws.Range("I99").Value = 1032.3334
$ws.Range("J99").Value = 197
$ws.Range("K99").Value = 3097.0002
$ws.Range("L99").Value = 591
$ws.Range("M99").Value = -1599.0002
$ws.Range("N99").Value = -3587

$ws.Range("H137").Value = 1161.7778
$ws.Range("I137").Value = 1072.6364
$ws.Range("J137").Value = 1301.8572
$ws.Range("K137").Value = 3217.9092
$ws.Range("L137").Value = 3905.5716
$ws.Range("M137").Value = -667.9092000000001
$ws.Range("N137").Value = -9005.571599999999

$ws.Range("H138").Value = 3331.257
$ws.Range("J138").Value = 3878.7693
$ws.Range("L138").Value = 11636.3079
$ws.Range("N138").Value = -21916.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2344.889
$ws.Range("I61").Value = 2263
$ws.Range("K61").Value = 2263
$ws.Range("M61").Value = -2051

$ws.Range("H80").Value = 44748.75
$ws.Range("J80").Value = 44748.75
$ws.Range("L80").Value = 44748.75
$ws.Range("N80").Value = -46744.75

$ws.Range("H83").Value = 44748.75
$ws.Range("J83").Value = 44748.75
$ws.Range("L83").Value = 134246.25
$ws.Range("N83").Value = -144230.25

$ws.Range("H97").Value = 2104.1333
$ws.Range("J97").Value = 6746
$ws.Range("L97").Value = 6746
$ws.Range("N97").Value = -7738

$ws.Range("H122").Value = 2579.3845
$ws.Range("I122").Value = 2290.818
$ws.Range("K122").Value = 6872.454000000001
$ws.Range("M122").Value = -4422.454000000001

$ws.Range("H136").Value = 2344.889
$ws.Range("I136").Value = 2263
$ws.Range("K136").Value = 6789
$ws.Range("M136").Value = -4239

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 802
$ws.Range("I5").Value = 802
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 802
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -689
$ws.Range("N5").ClearContents()

$ws.Range("H86").Value = 1887.2941
$ws.Range("I86").Value = 1741.7142
$ws.Range("J86").Value = 2566.6667
$ws.Range("K86").Value = 1741.7142
$ws.Range("L86").Value = 2566.6667
$ws.Range("M86").Value = -618.7141999999999
$ws.Range("N86").Value = -4812.6667

$ws.Range("H89").Value = 1887.2941
$ws.Range("I89").Value = 1741.7142
$ws.Range("J89").Value = 2566.6667
$ws.Range("K89").Value = 8708.571
$ws.Range("L89").Value = 12833.3335
$ws.Range("M89").Value = -3092.571
$ws.Range("N89").Value = -24065.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4616.6665
$ws.Range("I58").Value = 2166.6667
$ws.Range("J58").Value = 7066.6665
$ws.Range("K58").Value = 2166.6667
$ws.Range("L58").Value = 7066.6665
$ws.Range("M58").Value = -1963.6667
$ws.Range("N58").Value = -7472.6665

$ws.Range("H132").Value = 4627.2
$ws.Range("I132").Value = 3929.6667
$ws.Range("K132").Value = 11789.0001
$ws.Range("M132").Value = -9259.000100000001

$ws.Range("H134").Value = 2972.3333
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws.Range("H136").Value = 4616.6665
$ws.Range("I136").Value = 2166.6667
$ws.Range("J136").Value = 7066.6665
$ws.Range("K136").Value = 6500.000100000001
$ws.Range("L136").Value = 21199.9995
$ws.Range("M136").Value = -3950.000100000001
$ws.Range("N136").Value = -26299.9995

$ws.Range("H141").Value = 379999.5
$ws.Range("J141").Value = 379999.5
$ws.Range("L141").Value = 379999.5
$ws.Range("N141").Value = -390359.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 11645
$ws.Range("I80").Value = 4698.3335
$ws.Range("K80").Value = 14095.0005
$ws.Range("M80").Value = -13159.0005

$ws.Range("H83").Value = 11645
$ws.Range("I83").Value = 4698.3335
$ws.Range("K83").Value = 42285.0015
$ws.Range("M83").Value = -37605.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3409.6667
$ws.Range("I80").Value = 3498
$ws.Range("J80").Value = 3365.5
$ws.Range("K80").Value = 3498
$ws.Range("L80").Value = 3365.5
$ws.Range("M80").Value = -2500
$ws.Range("N80").Value = -5361.5

$ws.Range("H83").Value = 3409.6667
$ws.Range("I83").Value = 3498
$ws.Range("J83").Value = 3365.5
$ws.Range("K83").Value = 17490
$ws.Range("L83").Value = 16827.5
$ws.Range("M83").Value = -12498
$ws.Range("N83").Value = -26811.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2988.25
$ws.Range("I82").Value = 3000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2639

$ws.Range("H85").Value = 2988.25
$ws.Range("I85").Value = 3000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1752

$ws.Range("H136").Value = 3671.889
$ws.Range("I136").Value = 3664.8572
$ws.Range("J136").Value = 3696.5
$ws.Range("K136").Value = 10994.5716
$ws.Range("L136").Value = 11089.5
$ws.Range("M136").Value = -8444.571599999999
$ws.Range("N136").Value = -16189.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2567.25
$ws.Range("I81").Value = 2567.25
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5134.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -4073.5
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 2567.25
$ws.Range("I84").Value = 2567.25
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 25672.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -20368.5
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 1614.1428
$ws.Range("I122").Value = 1659.8
$ws.Range("K122").Value = 4979.4
$ws.Range("M122").Value = -2529.4

$ws.Range("H132").Value = 3388.1
$ws.Range("I132").Value = 3271.75
$ws.Range("J132").Value = 3465.6667
$ws.Range("K132").Value = 9815.25
$ws.Range("L132").Value = 10397.0001
$ws.Range("M132").Value = -7285.25
$ws.Range("N132").Value = -15457.0001

$ws.Range("H136").Value = 8606.429
$ws.Range("I136").Value = 9049
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 27147
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -24597
$ws.Range("N136").Value = -27600
